$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.982.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.437.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.27'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.53'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.527'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.434.21'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.17'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.30'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.875.58'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.001.29'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.430.18'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.89'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.06'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.26'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.09'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.83%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.53'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.31'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '614.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.556.89'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0940'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.82%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.96'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.88'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.29%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.42'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.98%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '149.88'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.23'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.21'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.75'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.79'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -8.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '142.04'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.60'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0521'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.596'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.46'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -7.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0231'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.98%  '
